# Update existing rows 8-15 (cols B-E) and append two new rows (16, 17)
# to the "lines_states" sheet, inserting line7 / line8 entries right after
# line6 (which pushes the former extr1..extr8 rows down by two positions).
#
# NOTE: the "name" column (B) previously held extr1..extr8 on rows 8-15; the
# new line7/line8 entries now occupy rows 8-9, and extr1..extr8 shift down to
# rows 10-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 (was extr1, now line7) ---
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# --- Row 9 (was extr2, now line8) ---
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# --- Row 10 (was extr3, now extr1) ---
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# --- Row 11 (was extr4, now extr2) ---
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# --- Row 12 (was extr5, now extr3) ---
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

# --- Row 13 (was extr6, now extr4) ---
$ws.Range("B13").Value = "extr4"
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

# --- Row 14 (was extr7, now extr5) ---
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

# --- Row 15 (was extr8, now extr6) ---
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# --- New row 16 (extr7; copy formatting from row 15, then set values) ---
$ws.Range("A15:E15").Copy($ws.Range("A16:E16"))
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

# --- New row 17 (extr8; copy formatting from row 16, then set values) ---
$ws.Range("A16:E16").Copy($ws.Range("A17:E17"))
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
